$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1937984496124031
$ws.Range("C2").Value = 0.5465116279069767
$ws.Range("J2").Value = 0.02325581395348837
$ws.Range("P2").Value = 0.1279069767441861
$ws.Range("S2").Value = 0.1085271317829457
$ws.Range("C3").Value = 0.01360544217687075
$ws.Range("J3").Value = 0.0272108843537415
$ws.Range("P3").Value = 0.782312925170068
$ws.Range("S3").Value = 0.1768707482993197
$ws.Range("J4").Value = 0.1025641025641026
$ws.Range("O4").Value = 0.02564102564102564
$ws.Range("P4").Value = 0.717948717948718
$ws.Range("B6").Value = 0.05339805825242718
$ws.Range("D6").Value = 0.009708737864077669
$ws.Range("F6").Value = 0.07766990291262135
$ws.Range("J6").Value = 0.3058252427184466
$ws.Range("O6").Value = 0.009708737864077669
$ws.Range("Q6").Value = 0.1116504854368932
$ws.Range("R6").Value = 0.0970873786407767
$ws.Range("S6").Value = 0.3349514563106796
$ws.Range("B7").Value = 0.1125827814569536
$ws.Range("D7").Value = 0.02649006622516556
$ws.Range("F7").Value = 0.0728476821192053
$ws.Range("J7").Value = 0.1456953642384106
$ws.Range("O7").Value = 0.02649006622516556
$ws.Range("Q7").Value = 0.1854304635761589
$ws.Range("R7").Value = 0.09933774834437085
$ws.Range("S7").Value = 0.3311258278145696
$ws.Range("B8").Value = 0.08500000000000001
$ws.Range("D8").Value = 0.01
$ws.Range("F8").Value = 0.0475
$ws.Range("J8").Value = 0.17
$ws.Range("O8").Value = 0.015
$ws.Range("Q8").Value = 0.1725
$ws.Range("R8").Value = 0.1275
$ws.Range("S8").Value = 0.3725
$ws.Range("B9").Value = 0.1090909090909091
$ws.Range("D9").Value = 0.01818181818181818
$ws.Range("E9").Value = 0.004545454545454545
$ws.Range("F9").Value = 0.02727272727272727
$ws.Range("J9").Value = 0.1318181818181818
$ws.Range("O9").Value = 0.004545454545454545
$ws.Range("Q9").Value = 0.2090909090909091
$ws.Range("R9").Value = 0.1
$ws.Range("S9").Value = 0.3954545454545454
$ws.Range("B10").Value = 0.09582309582309582
$ws.Range("D10").Value = 0.02211302211302211
$ws.Range("E10").Value = 0.0008190008190008191
$ws.Range("F10").Value = 0.06388206388206388
$ws.Range("J10").Value = 0.1506961506961507
$ws.Range("O10").Value = 0.01965601965601966
$ws.Range("Q10").Value = 0.2416052416052416
$ws.Range("R10").Value = 0.07698607698607698
$ws.Range("S10").Value = 0.3284193284193284
$ws.Range("G11").Value = 0.1164658634538153
$ws.Range("J11").Value = 0.1044176706827309
$ws.Range("K11").Value = 0.2048192771084337
$ws.Range("L11").Value = 0.5742971887550201
$ws.Range("G12").Value = 0.7019867549668874
$ws.Range("J12").Value = 0.1986754966887417
$ws.Range("K12").Value = 0.006622516556291391
$ws.Range("L12").Value = 0.06622516556291391
$ws.Range("S12").Value = 0.02649006622516556
$ws.Range("G13").Value = 0.5172413793103449
$ws.Range("J13").Value = 0.4482758620689655
$ws.Range("S13").Value = 0.03448275862068965
$ws.Range("F15").Value = 0.01415094339622642
$ws.Range("H15").Value = 0.160377358490566
$ws.Range("I15").Value = 0.08962264150943396
$ws.Range("J15").Value = 0.3584905660377358
$ws.Range("K15").Value = 0.0660377358490566
$ws.Range("M15").Value = 0.004716981132075472
$ws.Range("O15").Value = 0.07547169811320754
$ws.Range("S15").Value = 0.2311320754716981
$ws.Range("F16").Value = 0.04678362573099415
$ws.Range("H16").Value = 0.1929824561403509
$ws.Range("I16").Value = 0.07602339181286549
$ws.Range("J16").Value = 0.3684210526315789
$ws.Range("K16").Value = 0.1111111111111111
$ws.Range("M16").Value = 0.02923976608187134
$ws.Range("O16").Value = 0.06432748538011696
$ws.Range("S16").Value = 0.1111111111111111
$ws.Range("F17").Value = 0.01991150442477876
$ws.Range("H17").Value = 0.1615044247787611
$ws.Range("I17").Value = 0.084070796460177
$ws.Range("J17").Value = 0.4380530973451328
$ws.Range("K17").Value = 0.09513274336283185
$ws.Range("M17").Value = 0.01106194690265487
$ws.Range("N17").Value = 0.002212389380530973
$ws.Range("O17").Value = 0.05973451327433629
$ws.Range("S17").Value = 0.1283185840707965
$ws.Range("F18").Value = 0.01463414634146342
$ws.Range("H18").Value = 0.2146341463414634
$ws.Range("I18").Value = 0.1365853658536585
$ws.Range("J18").Value = 0.3902439024390244
$ws.Range("K18").Value = 0.07804878048780488
$ws.Range("M18").Value = 0.01951219512195122
$ws.Range("N18").Value = 0.004878048780487805
$ws.Range("O18").Value = 0.05853658536585366
$ws.Range("S18").Value = 0.08292682926829269
$ws.Range("F19").Value = 0.02420856610800745
$ws.Range("H19").Value = 0.1973929236499069
$ws.Range("I19").Value = 0.1070763500931099
$ws.Range("J19").Value = 0.3649906890130354
$ws.Range("K19").Value = 0.09869646182495345
$ws.Range("M19").Value = 0.01303538175046555
$ws.Range("N19").Value = 0.0009310986964618249
$ws.Range("O19").Value = 0.07914338919925512
$ws.Range("S19").Value = 0.1145251396648045
